# "Generate Report for Handback"
# Populates the handback columns (Latest Target File / Latest Handback File /
# Latest Handback DateTime) for the zh-cn and de-de worksheets, updates the
# status message on all three sheets from "Ready for handoff" to
# "Handed back: in sync with en-US", and widens a few columns to fit the
# newly-written content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Update status text wherever "Ready for handoff" currently appears ---
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsZhCn.Range("C2").Value     = $statusText
$wsDeDe.Range("C2").Value     = $statusText

# --- Hyperlink / file name color used for the custom "HyperLink" cell style ---
$hyperlinkColor = 15570276   # RGB(0x64,0x95,0xED) packed as BGR long

# --- zh-cn sheet: populate Latest Target File / Latest Handback File / DateTime ---
$wsZhCn.Range("I2").Value = "f7c5f93a-6302-4bf8-80fe-193effac2677.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ebe7936804450518b23f13753bacd6a2294c85b/e2e/f7c5f93a-6302-4bf8-80fe-193effac2677.md", "", "", "f7c5f93a-6302-4bf8-80fe-193effac2677.md")
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Color = $hyperlinkColor
$wsZhCn.Range("I2").Font.Name = "Calibri"

$wsZhCn.Range("J2").Value = "f7c5f93a-6302-4bf8-80fe-193effac2677.f6c6d10f2c2428ea0436459a440e31bdd6480fc0.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-19 19:06:59"

# --- de-de sheet: populate Latest Target File / Latest Handback File / DateTime ---
$wsDeDe.Range("I2").Value = "f7c5f93a-6302-4bf8-80fe-193effac2677.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ebe7936804450518b23f13753bacd6a2294c85b/e2e/f7c5f93a-6302-4bf8-80fe-193effac2677.md", "", "", "f7c5f93a-6302-4bf8-80fe-193effac2677.md")
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Color = $hyperlinkColor
$wsDeDe.Range("I2").Font.Name = "Calibri"

$wsDeDe.Range("J2").Value = "f7c5f93a-6302-4bf8-80fe-193effac2677.f6c6d10f2c2428ea0436459a440e31bdd6480fc0.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-19 19:07:13"

# --- Widen columns to fit the newly populated / lengthened content ---
# (ColumnWidth is internally quantized to 1/6-character increments, so the
# inputs below are chosen to land as close as possible to the true target
# widths of ~29.98, ~39.37 and 40 characters.)
$wsOverview.Columns.Item(5).ColumnWidth  = 29.166666666666668   # E
$wsOverview.Columns.Item(6).ColumnWidth  = 29.166666666666668   # F

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.166666666666668   # C
$wsZhCn.Columns.Item(9).ColumnWidth  = 38.5                 # I
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664   # J

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.166666666666668   # C
$wsDeDe.Columns.Item(9).ColumnWidth  = 38.5                 # I
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664   # J
